# Refresh the Ccl3-Ccr5 LR-pair sheet with updated TPM-derived NATMI output.
# Existing rows 2-4 (ECs->ECs, ECs->FAPs, ECs->MuSCs) are replaced with the
# new 4-combination matrix (ECs/MuSCs sending -> ECs/FAPs target), and a new
# row 5 is added to hold the 4th combination.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ccl3"
$ws.Range("C2").Value = "Ccr5"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.03315566666666667
$ws.Range("H2").Value = 0.099467
$ws.Range("I2").Value = 0.1557603470145164
$ws.Range("J2").Value = 0.1557603470145164
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.001937666666666667
$ws.Range("N2").Value = 0.005813
$ws.Range("O2").Value = 0.0230007399171451
$ws.Range("P2").Value = 0.02300073991714511
$ws.Range("Q2").Value = 0.00006424463011111111
$ws.Range("R2").Value = 0.000578201671
$ws.Range("S2").Value = 0.003582603231085159
$ws.Range("T2").Value = 0.00358260323108516

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ccl3"
$ws.Range("C3").Value = "Ccr5"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.03315566666666667
$ws.Range("H3").Value = 0.099467
$ws.Range("I3").Value = 0.1557603470145164
$ws.Range("J3").Value = 0.1557603470145164
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.082306
$ws.Range("N3").Value = 0.246918
$ws.Range("O3").Value = 0.9769992600828549
$ws.Range("P3").Value = 0.976999260082855
$ws.Range("Q3").Value = 0.002728910300666667
$ws.Range("R3").Value = 0.024560192706
$ws.Range("S3").Value = 0.1521777437834312
$ws.Range("T3").Value = 0.1521777437834312

$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Ccl3"
$ws.Range("C4").Value = "Ccr5"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1797076666666667
$ws.Range("H4").Value = 0.539123
$ws.Range("I4").Value = 0.8442396529854836
$ws.Range("J4").Value = 0.8442396529854836
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.001937666666666667
$ws.Range("N4").Value = 0.005813
$ws.Range("O4").Value = 0.0230007399171451
$ws.Range("P4").Value = 0.02300073991714511
$ws.Range("Q4").Value = 0.0003482135554444444
$ws.Range("R4").Value = 0.003133921999
$ws.Range("S4").Value = 0.01941813668605994
$ws.Range("T4").Value = 0.01941813668605995

$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Ccl3"
$ws.Range("C5").Value = "Ccr5"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.1797076666666667
$ws.Range("H5").Value = 0.539123
$ws.Range("I5").Value = 0.8442396529854836
$ws.Range("J5").Value = 0.8442396529854836
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.082306
$ws.Range("N5").Value = 0.246918
$ws.Range("O5").Value = 0.9769992600828549
$ws.Range("P5").Value = 0.976999260082855
$ws.Range("Q5").Value = 0.01479101921266667
$ws.Range("R5").Value = 0.133119172914
$ws.Range("S5").Value = 0.8248215162994237
$ws.Range("T5").Value = 0.8248215162994238

"Updated range A1:T5 -> " + $ws.Range("A1:T5").Address()
